# Generate Report for Handoff
#
# Refresh the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamp for the
# 2af883b4-c805-4bd7-af30-384957281dcd file (row 6 on every sheet) to reflect a
# freshly-generated handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G), row 6.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-09-02 14:52:02"

# zh-cn sheet: "Latest Handoff Datetime" column (H), row 6.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-09-02 14:51:56"

# de-de sheet: "Latest Handoff Datetime" column (H), row 6.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-09-02 14:52:02"
